# Append new experiment log rows (2016-04-26 runs, n_iterator 1000/2000/3000)
# to the "logs" sheet - update result from 001 - 004.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$features = "12 features: %ascii-adp, %digit-adp, digit-adp/ascii-adp, %keyword-name, %keyword-address, %keyword-phone, bfirst-character-digit, bfirst-character-ascii, blast-character-digit, blast-character-ascii, b#ascii >= 6, b#digit >= 7"
$nn = "Neural-Network"

# Each entry: Time, Model description (layers/lr/rule/iterations), ClassifyAcc, SegmentAcc
$newRows = @(
    ,@("20160426_092146", "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000", 0.914191419141914, 0.47)
    ,@("20160426_093615", "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000", 0.914191419141914, 0.49)
    ,@("20160426_095136", "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000", 0.914191419141914, 0.47)
    ,@("20160426_100639", "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000", 0.914191419141914, 0.46)
    ,@("20160426_102150", "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000", 0.914191419141914, 0.53)
    ,@("20160426_110903", "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000", 0.907590759075908, 0.46)
    ,@("20160426_113732", "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000", 0.914191419141914, 0.47)
    ,@("20160426_120658", "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000", 0.914191419141914, 0.47)
    ,@("20160426_123554", "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000", 0.910891089108911, 0.47)
    ,@("20160426_130510", "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000", 0.910891089108911, 0.47)
    ,@("20160426_134439", "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 3000", 0.894389438943894, 0.44)
    ,@("20160426_142957", "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 3000", 0.907590759075908, 0.46)
    ,@("20160426_151439", "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 3000", 0.900990099009901, 0.45)
    ,@("20160426_155825", "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 3000", 0.897689768976898, 0.44)
    ,@("20160426_164244", "2 layers: [10-Sigmoid, 2-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 3000", 0.900990099009901, 0.44)
)

$startRow = 12
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $time = $row[0]
    $model = $row[1]
    $classifyAcc = $row[2]
    $segmentAcc = $row[3]

    $ws.Cells.Item($r, 1).Value = $time          # Time
    $ws.Cells.Item($r, 2).Value = $features       # NameFeatures
    $ws.Cells.Item($r, 3).Value = $features       # AddressFeatures
    $ws.Cells.Item($r, 4).Value = $features       # PhoneFeatures
    $ws.Cells.Item($r, 5).Value = $nn             # NameModelType
    $ws.Cells.Item($r, 6).Value = $model          # NameModel
    $ws.Cells.Item($r, 7).Value = $nn             # AddressModelType
    $ws.Cells.Item($r, 8).Value = $model          # AddressModel
    $ws.Cells.Item($r, 9).Value = $nn             # PhoneModelType
    $ws.Cells.Item($r, 10).Value = $model         # PhoneModel
    $ws.Cells.Item($r, 11).Value = $classifyAcc   # ClassifyAcc
    $ws.Cells.Item($r, 12).Value = $segmentAcc    # SegmentAcc
}

